# Mise à jour de certains champs de Modules et de Professeurs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "Matières enseignés" column (adds shared string + extends row/dimension)
$ws.Range("E1").Value = "Matières enseignés"

# Explicit column widths for C, D and E as set in the updated sheet
$ws.Columns.Item(3).ColumnWidth = 26.666666666666668
$ws.Columns.Item(4).ColumnWidth = 14.833333333333334
$ws.Columns.Item(5).ColumnWidth = 30.833333333333332

# Update the active selection/view to match the saved workbook state
$ws.Range("E6").Select()
